$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 height: 19.5 -> 20.25 (to match the rest of the rows' post-edit height)
$ws.Rows.Item(2).RowHeight = 20.25

# Row 3 ("News page implementaions"): status flips to Completed and start/end
# dates are filled in with 20-07-2025
$ws.Range("D3").Value = "Completed"
$ws.Range("E3").Value = "20-07-2025"
$ws.Range("F3").Value = "20-07-2025"

# Row 4 ("Tools"): no text change, just drop the cell border so it matches
# the rest of the "Not Started" cells that have no border
$ws.Cells.Item(4,4).Borders.LineStyle = -4142

# Row 7 ("Snippet implementaions"): status flips to Completed and start/end
# dates are filled in with 20-07-2025
$ws.Range("D7").Value = "Completed"
$ws.Range("E7").Value = "20-07-2025"
$ws.Range("F7").Value = "20-07-2025"

# Row 13 ("Messages fine tune for all form"): end/start date corrected from
# 18-07-2025 to 19-07-2025
$ws.Range("E13").Value = "19-07-2025"
$ws.Range("F13").Value = "19-07-2025"

# New row 14: "Env for development server" task, High priority, Not Started
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Env for development server"
$ws.Range("C14").Value = "High"
$ws.Range("D14").Value = "Not Started"
